# Auto-generated edit script to update Kujata_Profits sheet values
# per scheduled runner profit recalculation
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2371.7144
$ws.Range("I18").Value = 2675
$ws.Range("J18").Value = 1967.3334
$ws.Range("K18").Value = 2675
$ws.Range("L18").Value = 1967.3334
$ws.Range("M18").Value = -2391
$ws.Range("N18").Value = -2535.3334

$ws.Range("H92").Value = 1558.9445
$ws.Range("I92").Value = 1614
$ws.Range("K92").Value = 1614
$ws.Range("M92").Value = -366

$ws.Range("H111").Value = 3585.8
$ws.Range("I111").Value = 3585.8
$ws.Range("K111").Value = 10757.4
$ws.Range("M111").Value = -7690.400000000001

$ws.Range("H138").Value = 513083.97
$ws.Range("I138").Value = 1324.2222
$ws.Range("J138").Value = 751316.9399999999
$ws.Range("K138").Value = 3972.6666
$ws.Range("L138").Value = 2253950.82
$ws.Range("M138").Value = 1167.3334
$ws.Range("N138").Value = -2264230.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3821.8
$ws.Range("I32").Value = 3328.3281
$ws.Range("K32").Value = 3328.3281
$ws.Range("M32").Value = -3041.3281

$ws.Range("H74").Value = 1455.725
$ws.Range("I74").Value = 812.86365
$ws.Range("J74").Value = 2241.4443
$ws.Range("K74").Value = 812.86365
$ws.Range("L74").Value = 2241.4443
$ws.Range("M74").Value = 61.13634999999999
$ws.Range("N74").Value = -3989.4443

$ws.Range("H77").Value = 1455.725
$ws.Range("I77").Value = 812.86365
$ws.Range("J77").Value = 2241.4443
$ws.Range("K77").Value = 4064.31825
$ws.Range("L77").Value = 11207.2215
$ws.Range("M77").Value = 303.6817499999997
$ws.Range("N77").Value = -19943.2215

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 62501276
$ws.Range("I99").Value = 83334536
$ws.Range("K99").Value = 83334536
$ws.Range("M99").Value = -83333038

$ws.Range("H138").Value = 67440
$ws.Range("J138").Value = 67440
$ws.Range("L138").Value = 67440
$ws.Range("N138").Value = -77720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 6203
$ws.Range("I17").Value = 4300
$ws.Range("K17").Value = 4300
$ws.Range("M17").Value = -4126

$ws.Range("H25").Value = 8812.6
$ws.Range("I25").Value = 4011
$ws.Range("K25").Value = 4011
$ws.Range("M25").Value = -3837

$ws.Range("H41").Value = 6960
$ws.Range("I41").Value = 2575
$ws.Range("J41").Value = 24500
$ws.Range("K41").Value = 2575
$ws.Range("L41").Value = 24500
$ws.Range("M41").Value = -2147
$ws.Range("N41").Value = -25356

$ws.Range("H50").Value = 26892
$ws.Range("J50").Value = 26892
$ws.Range("L50").Value = 26892
$ws.Range("N50").Value = -28142

$ws.Range("H51").Value = 16666.666
$ws.Range("J51").Value = 22500
$ws.Range("L51").Value = 22500
$ws.Range("N51").Value = -23972

$ws.Range("H60").Value = 3150
$ws.Range("I60").Value = 3150
$ws.Range("K60").Value = 3150
$ws.Range("M60").Value = -2639

$ws.Range("H61").Value = 16666.666
$ws.Range("J61").Value = 22500
$ws.Range("L61").Value = 22500
$ws.Range("N61").Value = -23196

$ws.Range("H64").Value = 33000
$ws.Range("J64").Value = 33000
$ws.Range("L64").Value = 33000
$ws.Range("N64").Value = -33496

$ws.Range("H67").Value = 33000
$ws.Range("J67").Value = 33000
$ws.Range("L67").Value = 33000
$ws.Range("N67").Value = -34716

$ws.Range("H135").Value = 34698
$ws.Range("J135").Value = 34698
$ws.Range("L135").Value = 34698
$ws.Range("N135").Value = -44838

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19233560
$ws.Range("J131").Value = 3559.85
$ws.Range("L131").Value = 10679.55
$ws.Range("N131").Value = -20759.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3238
$ws.Range("I80").Value = 1796.6666
$ws.Range("J80").Value = 5400
$ws.Range("K80").Value = 1796.6666
$ws.Range("L80").Value = 5400
$ws.Range("M80").Value = -798.6666
$ws.Range("N80").Value = -7396

$ws.Range("H83").Value = 3238
$ws.Range("I83").Value = 1796.6666
$ws.Range("J83").Value = 5400
$ws.Range("K83").Value = 8983.333000000001
$ws.Range("L83").Value = 27000
$ws.Range("M83").Value = -3991.333000000001
$ws.Range("N83").Value = -36984

$ws.Range("H107").Value = 861.2778
$ws.Range("J107").Value = 783.1667
$ws.Range("L107").Value = 783.1667
$ws.Range("N107").Value = -4623.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 801.94116
$ws.Range("I16").Value = 664.4375
$ws.Range("J16").Value = 3002
$ws.Range("K16").Value = 664.4375
$ws.Range("L16").Value = 3002
$ws.Range("M16").Value = -494.4375
$ws.Range("N16").Value = -3342

$ws.Range("H55").Value = 368.17648
$ws.Range("I55").Value = 269.45456
$ws.Range("K55").Value = 269.45456
$ws.Range("M55").Value = -96.45456000000001

$ws.Range("H64").Value = 18448.5
$ws.Range("J64").Value = 18448.5
$ws.Range("L64").Value = 18448.5
$ws.Range("N64").Value = -18898.5

$ws.Range("H67").Value = 18448.5
$ws.Range("J67").Value = 18448.5
$ws.Range("L67").Value = 18448.5
$ws.Range("N67").Value = -20008.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 60000
$ws.Range("J63").Value = 60000
$ws.Range("L63").Value = 60000
$ws.Range("N63").Value = -61248

$ws.Range("H66").Value = 60000
$ws.Range("J66").Value = 60000
$ws.Range("L66").Value = 180000
$ws.Range("N66").Value = -186240

$ws.Range("H100").Value = 838
$ws.Range("I100").Value = 1067.8572
$ws.Range("J100").Value = 569.8333
$ws.Range("K100").Value = 2135.7144
$ws.Range("L100").Value = 1139.6666
$ws.Range("M100").Value = -1594.7144
$ws.Range("N100").Value = -2221.6666

$ws.Range("H122").Value = 17858540
$ws.Range("I122").Value = 17858540
$ws.Range("K122").Value = 53575620
$ws.Range("M122").Value = -53573170

$ws.Range("H128").Value = 74857.5
$ws.Range("J128").Value = 74857.5
$ws.Range("L128").Value = 74857.5
$ws.Range("N128").Value = -84817.5
